$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header question text (strip wrapping quotes / trailing comma) ---
$ws.Range("C1").Value = 'Does society_name encompasses community sites? Respond one word (''yes'' or ''no'') only plus provide a justification for the answer also after a comma.'
$ws.Range("D1").Value = 'Is society_name influential on state or local policy? Respond one word (''yes'' or ''no'') only plus provide a justification for the answer also after a comma.'
$ws.Range("E1").Value = 'Does society_name provide engagement opportunity with leadership? Respond one word (''yes'' or ''no'') only plus provide a justification for the answer also after a comma.'
$ws.Range("F1").Value = 'Does society_name provide support for clinical trial recruitment? Respond one word (''yes'' or ''no'') only plus provide a justification for the answer also after a comma.'
$ws.Range("G1").Value = 'Does society_name provide engagement opportunity with payors? Respond one word (''yes'' or ''no'') only plus provide a justification for the answer also after a comma.'
$ws.Range("H1").Value = 'Does society_name include area experts on its board? Respond one word (''yes'' or ''no'') only plus provide a justification for the answer also after a comma.'
$ws.Range("I1").Value = 'Is society_name involved in therapeutic research collaborations? Respond one word (''yes'' or ''no'') only plus provide a justification for the answer also after a comma.'
$ws.Range("J1").Value = 'Does society_name include top therapeutic area experts on its board? Respond with one word (''yes'' or ''no'') only plus provide a justification for the answer also after a comma.'
$ws.Range("K1").Value = 'Name the Region where the society_name is from? Just name the US state or the Region in one word for the answer.'

# --- Row 2 ---
$ws.Range("A2").Value = 'Soleo Health'
$ws.Range("C2").Value = 'No, Soleo Health does not encompass community sites. Soleo Health is a specialty pharmacy and infusion services provider, focusing on clinical expertise and patient care rather than community-driven initiatives.'
$ws.Range("D2").Value = 'No, Soleo Health is primarily focused on healthcare services and does not typically engage in policy advocacy or influence at the state or local level.'
$ws.Range("E2").Value = 'Yes, Soleo Health provides engagement opportunity with leadership. Soleo Health encourages collaboration between leadership and employees, offering opportunities for interaction and involvement in decision-making processes.'
$ws.Range("F2").Value = 'No, Soleo Health does not provide support for clinical trial recruitment. Soleo Health specializes in specialty infusion services, not clinical trial recruitment.'
$ws.Range("G2").Value = 'No, Soleo Health does not provide engagement opportunities with payors. They focus on providing specialty infusion services.'
$ws.Range("H2").Value = 'No, justification: Soleo Health does not publicly disclose information about area experts on its board.'
$ws.Range("I2").Value = 'No, Soleo Health is not involved in therapeutic research collaborations. Soleo Health focuses on providing specialty pharmacy and infusion services to patients.'
$ws.Range("J2").Value = 'No, justification: The information about specific board members of Soleo Health is not publicly available to confirm if top therapeutic area experts are included.'
$ws.Range("K2").Value = 'Texas'

# --- Row 3 ---
$ws.Range("A3").Value = 'University of Miami'
$ws.Range("C3").Value = 'Yes, University of Miami encompasses community sites. The university is an integral part of the local community and engages with various community initiatives and programs.'
$ws.Range("D3").Value = 'No, The University of Miami is not influential on state or local policy. The university focuses more on education and research rather than policy advocacy.'
$ws.Range("E3").Value = 'yes, The University of Miami provides engagement opportunities with leadership through various campus leadership programs, student organizations, and workshops, allowing students to develop their leadership skills and engage with leaders in different fields.'
$ws.Range("F3").Value = 'No, The University of Miami does not support clinical trial recruitment. The institution focuses more on conducting research and trials rather than recruiting participants.'
$ws.Range("G3").Value = 'No, the University of Miami does not provide engagement opportunity with payors. The society primarily focuses on academic research and student programs.'
$ws.Range("H3").Value = 'No, the University of Miami society does not have area experts on its board. The board primarily consists of university administrators, faculty members, and external community leaders.'
$ws.Range("I3").Value = 'No, there is no available information indicating that the society ''University of Miami'' is involved in therapeutic research collaborations.'
$ws.Range("J3").Value = 'No, the University of Miami does not include top therapeutic area experts on its board. The board of directors typically consists of individuals with expertise in various fields like education, business, and governance, rather than focusing solely on therapeutic areas.'
$ws.Range("K3").Value = 'Florida'

# --- Row 4 ---
$ws.Range("A4").Value = 'Dava Oncology, LP'
$ws.Range("C4").Value = 'No, justification: Oncology-focused society focusing on medical professionals and researchers, not community sites.'
$ws.Range("D4").Value = 'No, Dava Oncology, LP is not influential on state or local policy. Dava Oncology, LP is a private healthcare company focused on providing oncology products and services, and typically does not have direct influence on state or local policy decisions.'
$ws.Range("E4").Value = 'No, justification: Dava Oncology, LP does not provide public information about engagement opportunities with their leadership.'
$ws.Range("F4").Value = 'No, Dava Oncology, LP does not provide support for clinical trial recruitment. Justification: Dava Oncology, LP is a pharmaceutical company that focuses on oncology medications rather than clinical trial recruitment services.'
$ws.Range("G4").Value = 'No, Dava Oncology, LP does not provide engagement opportunity with payors. Justification: Specializes in pharmaceutical industry.'
$ws.Range("H4").Value = 'No, no area experts on the board. Dava Oncology, LP primarily focuses on oncology expertise.'
$ws.Range("I4").Value = 'Yes, Dava Oncology, LP is involved in therapeutic research collaborations. Dava Oncology, LP often collaborates with other organizations to advance oncology research and develop new therapies.'
$ws.Range("J4").Value = 'No, Dava Oncology, LP does not include top therapeutic area experts on its board. , There is no public information available to suggest that top therapeutic area experts are part of the board of Dava Oncology, LP.'
$ws.Range("K4").Value = 'Texas'

# --- Membership counts (numeric-looking text, keep as text like the original) ---
$ws.Range("B2").Formula = "'1500"
$ws.Range("B3").Formula = "'13"
$ws.Range("B4").Formula = "'600"
$ws.Range("B2:B4").Style = "Normal"
